# Generate Report for Handoff
# - Update status text from "Handed back: in sync with en-US" to "Ready for handoff"
# - Refresh the handoff/generate timestamps
# - Narrow the "status" columns (they had been the same width as timestamp columns)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
# Overview!G2 and de-de!H2 shared the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" stamp
$wsOverview.Range("G2").Value = "2016-09-09 13:05:23"
$wsDeDe.Range("H2").Value = "2016-09-09 13:05:23"

# zh-cn!H2 "Latest Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-09-09 13:05:00"

# --- Column widths: status columns narrowed from 29.9777050018311 to 17.2159881591797 ---
# (ColumnWidth is quantized internally to 1/6-character steps; 16.25 is the input
# that lands on the closest representable width to the target.)
$wsOverview.Range("E1").ColumnWidth = 16.25
$wsOverview.Range("F1").ColumnWidth = 16.25
$wsZhCn.Range("C1").ColumnWidth = 16.25
$wsDeDe.Range("C1").ColumnWidth = 16.25
